# Mindforge_Student_Experience_Project_Tracking.xlsx update
# "Added Backend and new docs"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sprint & Task Tracking")
$ws2 = $wb.Worksheets.Item("Sprint Summary")

# --- Sheet 1: Sprint & Task Tracking ---

# Date-like text cells must stay plain text (not auto-converted to a date
# serial number by Excel's smart entry), so force text format first.
$dateCells = @("L3", "M2", "M3", "L4", "M4", "L5", "M5")
foreach ($addr in $dateCells) {
    $ws1.Range($addr).NumberFormat = "@"
}

# Row 2 (Task 1.1): completed date moves a day, tracker notes rewritten for NestJS rebuild
$ws1.Range("M2").Value = "2026-02-11"
$ws1.Range("P2").Value = "REBUILT on NestJS+TypeScript (locked stack). Modular: /modules/{auth,student,attendance,activities} with controller/service/repository/policy pattern. Global: Helmet; HTTPS enforce; CORS; JSON-only; class-validator DTO pipes; GlobalExceptionFilter; RequestId; Logging; ThrottlerGuard; AuthGuard; AuthorizationGuard; AuditService. Health: GET /health. Auth: POST /v1/auth/* with DTO validation. OpenAPI/Swagger at /api/docs. Zero TS build errors."

# Row 3 (Task 1.2): now Done, with start/completed dates and tracker notes
$ws1.Range("K3").Value = "Done"
$ws1.Range("L3").Value = "2026-02-11"
$ws1.Range("M3").Value = "2026-02-11"
$ws1.Range("P3").Value = "POST /v1/auth/mpin/verify: bcrypt MPIN hash via AuthRepository (in-memory; DB Task 2.1). Success: JWT (1h, studentId). Failure: 401 generic. Lockout: 5 attempts -> 15 min (403). Rate limit: 10/60s. Audit: LOGIN_SUCCESS/FAILURE/LOCKOUT. No MPIN in logs/response."

# Row 4 (Task 1.3): now Done, with start/completed dates and tracker notes
$ws1.Range("K4").Value = "Done"
$ws1.Range("L4").Value = "2026-02-11"
$ws1.Range("M4").Value = "2026-02-11"
$ws1.Range("P4").Value = "POST /v1/auth/lockout/status: fully implemented — returns isLocked, lockedUntil, attemptsRemaining, maxAttempts (5), lockoutDurationMinutes (15). POST /v1/auth/forgot-mpin: returns 202 Accepted — v1 entry point; full OTP deferred. Audit logged. Error shape consistent."

# Row 5 (Task 1.4): now Done, with start/completed dates and tracker notes
$ws1.Range("K5").Value = "Done"
$ws1.Range("L5").Value = "2026-02-11"
$ws1.Range("M5").Value = "2026-02-11"
$ws1.Range("P5").Value = "AuthGuard verifies JWT via JwtService.verifyAsync. Expired → 401 TOKEN_EXPIRED; invalid → 401 INVALID_TOKEN; missing → 401 UNAUTHORIZED. Extracts studentId from JWT payload → request.student. @Student() decorator. @Public() bypass. GET /v1/student/me protected demo endpoint. No cross-student access."

# --- Sheet 2: Sprint Summary ---
# Sprint 1 planned SP drops 8 -> 4, Done tasks count 1 -> 2 (Remaining recalculates via formula)
$ws2.Range("D2").Value = 4
$ws2.Range("F2").Value = 2
